$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 42; existing rows 42-69 shift down to 43-70.
$ws.Rows.Item(42).Insert()

# Populate the new row 42 with the "new" data point.
$ws.Range("A42").Value = 9
$ws.Range("B42").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C42").Value = "Metropolitana"
$ws.Range("D42").Value = 44488
$ws.Range("D42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = 100112022
$ws.Range("G42").Value = "Arveja Verde"
$ws.Range("H42").Value = "Perfection"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 25
$ws.Range("K42").Value = 24000
$ws.Range("L42").Value = 25000
$ws.Range("M42").Value = 24480
$ws.Range("N42").Value = "`$/malla 25 kilos"
$ws.Range("O42").Value = "Provincia de Huasco"
$ws.Range("P42").Value = 979
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"
